$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '25.927.26'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  +0.25%  '

# Row 3
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.640.53'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  +0.09%  '

# Row 4
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '1.003'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  +0.22%  '

# Row 5
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '215.15'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.08%  '

# Row 6
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.5056'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +0.51%  '

# Row 7
$ws.Range('E7').Value = '  +0.02%  '

# Row 8
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.2561'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -0.50%  '

# Row 9
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.06368'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -0.06%  '

# Row 10
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '19.43'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -0.22%  '

# Row 11
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.07756'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +0.23%  '

# Row 12
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '1.653.73'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +0.09%  '

# Row 13
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '4.276'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +0.45%  '

# Row 14
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.5438'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -0.21%  '

# Row 15
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.0₅7805'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -1.02%  '

# Row 16
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '64.24'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +0.13%  '

# Row 17
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '25.972.23'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +0.23%  '

# Row 18
$ws.Range('E18').Value = '  -0.02%  '

# Row 19
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '196.92'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -2.45%  '

# Row 20
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '4.430'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +1.05%  '

# Row 21
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '9.933'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +0.53%  '

# Row 22
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '6.034'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +1.09%  '

# Row 23
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '1.004'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -0.02%  '

# Row 24
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '1.894'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +1.59%  '

# Row 25
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '140.64'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -0.18%  '

# Row 26
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '0.1166'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +2.74%  '

# Row 27
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '6.876'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +1.53%  '

# Row 28
$ws.Range('E28').Value = '  +0.21%  '

# Row 29
$ws.Range('E29').Value = '  -0.56%  '

# Row 30
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '0.04949'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -0.37%  '

# Row 31
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '3.254'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -0.47%  '

# Row 32
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '3.180'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -0.45%  '

# Row 33
$ws.Range('E33').Value = '  -0.71%  '

# Row 34
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '2.362'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -0.34%  '

# Row 35
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.8936'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +0.35%  '

# Row 36
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '2.590'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -1.43%  '

# Row 37
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '1.133.43'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -1.36%  '

# Row 38
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.5437'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -2.75%  '

# Row 39
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.01557'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -0.42%  '

# Row 40
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '1.003'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -0.12%  '

# Row 41
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '2.542'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -1.11%  '

# Row 42
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.8189'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +1.65%  '

# Row 43
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '5.579'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -1.82%  '

# Row 44
$ws.Range('B44').Value = 'BabyDogeCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.0₈127'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +7.35%  '

# Row 45
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '99.54'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -0.24%  '

# Row 46
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '1.777.51'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -0.07%  '

# Row 47
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.4537'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +0.21%  '

# Row 48
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '1.001'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -0.81%  '

# Row 49
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '54.60'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -0.24%  '

# Row 50
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.05071'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +0.41%  '

# Row 51
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '1.005'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +0.35%  '
